$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
# B1 keeps its existing label ("Algorithm"); C1:L1 become the ten
# mean/std headers (two per base period instead of one).
$headers = @(
    "One Year Base mean", "One Year Base std",
    "Two Year Base mean", "Two Year Base std",
    "Three Year Base mean", "Three Year Base std",
    "Five Year Base mean", "Five Year Base std",
    "Ten Year Base mean", "Ten Year Base std"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 3 + $i   # C = 3
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# The newly introduced H1:L1 cells have no prior formatting (they didn't
# exist before), so copy the bold/border/center header style from the
# existing G1 header cell onto them.
$ws.Range("G1").Copy()
$ws.Range("H1:L1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows (rows 2-8) ----------------------------------------------
# Algorithm names: CART -> DTREE, NB dropped entirely (row count 8 -> 7).
$algorithms = @("LR", "LDA", "KNN", "DTREE", "RTREE", "XTREE", "SVM")

$values = @(
    @(0.9095535075871932, 0.01015304968933212, 0.8988380147739339, 0.007445084132257713, 0.8880148310533601, 0.01197825427050512, 0.8786682992495848, 0.01919816365316406, 0.8744435018663047, 0.01341383347521556),
    @(0.9119935912714174, 0.00876623223729628,  0.9055813900670563, 0.006210177132949237, 0.8929142223800717, 0.01129297585710078, 0.8862379704416566, 0.0168819316220737,   0.8763426083022283, 0.01390509824264928),
    @(0.8579825746140198, 0.01101777327952865,  0.8521360314783587, 0.01647407510778245,  0.8605412538424797, 0.01671472183098274, 0.8720452397319127, 0.01491086224952488,  0.8658771632168307, 0.00999786833719731),
    @(0.8781584174147931, 0.01617647336543993,  0.8748948166317476, 0.01276964486882566,  0.8771631170930645, 0.01420545074268068, 0.8691975282121784, 0.01358393351739461,  0.8713460015835313, 0.008888284668409033),
    @(0.8636791928179868, 0.01454898135759998,  0.8750648701744824, 0.0150767445012826,   0.8670205013900285, 0.01519433324703699, 0.8591707481239619, 0.009401642012249946, 0.8539910643592353, 0.01490232695084443),
    @(0.9100389290537857, 0.00703443734668265,  0.9042326014501395, 0.01091795644045563,  0.8953605501328795, 0.01365842593497318, 0.8960853812224322, 0.01582705342689686,  0.8872825472231648, 0.01530033115993424),
    @(0.9105280580493101, 0.007528987134322255, 0.9038964689049005, 0.008643830714888164, 0.900786254026184,  0.01389135375833084, 0.8949483015409292, 0.01819129910625266,  0.892991177468612,  0.01224487369078438)
)

# Remove the extra (8th) data row - NB row was dropped, shrinking the
# table from 8 data rows to 7.
$ws.Rows.Item(9).Delete()

for ($r = 0; $r -lt 7; $r++) {
    $row = 2 + $r
    $ws.Cells.Item($row, 1).Value = $r              # column A: 0-based index
    $ws.Cells.Item($row, 2).Value = $algorithms[$r]  # column B: algorithm name

    for ($c = 0; $c -lt 10; $c++) {
        $ws.Cells.Item($row, 3 + $c).Value = $values[$r][$c]
    }
}
